$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = [double]"16.37389066666667"
    "H2" = [double]"49.121672"
    "I2" = [double]"0.09466313117816218"
    "J2" = [double]"0.09466313117816218"
    "M2" = [double]"27.681071"
    "N2" = [double]"83.04321300000001"
    "O2" = [double]"0.05045805550111082"
    "P2" = [double]"0.05045805550111081"
    "Q2" = [double]"453.2468300902374"
    "R2" = [double]"4079.221470812137"
    "S2" = [double]"0.004776517526896641"
    "T2" = [double]"0.004776517526896641"
    "G3" = [double]"16.37389066666667"
    "H3" = [double]"49.121672"
    "I3" = [double]"0.09466313117816218"
    "J3" = [double]"0.09466313117816218"
    "O3" = [double]"0.0001771869602491167"
    "P3" = [double]"0.0001771869602491166"
    "Q3" = [double]"1.591607668362667"
    "R3" = [double]"14.324469015264"
    "S3" = [double]"1.677307246112194E-05"
    "T3" = [double]"1.677307246112193E-05"
    "G4" = [double]"16.37389066666667"
    "H4" = [double]"49.121672"
    "I4" = [double]"0.09466313117816218"
    "J4" = [double]"0.09466313117816218"
    "M4" = [double]"272.2666776666667"
    "N4" = [double]"816.800033"
    "O4" = [double]"0.4962975288350554"
    "P4" = [double]"0.4962975288350553"
    "Q4" = [double]"4458.064812290576"
    "R4" = [double]"40122.58331061518"
    "S4" = [double]"0.04698107807551057"
    "T4" = [double]"0.04698107807551057"
    "G5" = [double]"16.37389066666667"
    "H5" = [double]"49.121672"
    "I5" = [double]"0.09466313117816218"
    "J5" = [double]"0.09466313117816218"
    "M5" = [double]"11.73516533333333"
    "N5" = [double]"35.205496"
    "O5" = [double]"0.02139128300722342"
    "P5" = [double]"0.02139128300722341"
    "Q5" = [double]"192.1503141232569"
    "R5" = [double]"1729.352827109312"
    "S5" = [double]"0.002024965829381982"
    "T5" = [double]"0.002024965829381982"
    "G6" = [double]"16.37389066666667"
    "H6" = [double]"49.121672"
    "I6" = [double]"0.09466313117816218"
    "J6" = [double]"0.09466313117816218"
    "M6" = [double]"236.8155566666667"
    "N6" = [double]"710.44667"
    "O6" = [double]"0.4316759456963613"
    "P6" = [double]"0.4316759456963613"
    "Q6" = [double]"3877.592033025805"
    "R6" = [double]"34898.32829723224"
    "S6" = [double]"0.04086379667391186"
    "T6" = [double]"0.04086379667391186"
    "I7" = [double]"0.4193879037829277"
    "J7" = [double]"0.4193879037829278"
    "M7" = [double]"27.681071"
    "N7" = [double]"83.04321300000001"
    "O7" = [double]"0.05045805550111082"
    "P7" = [double]"0.05045805550111081"
    "Q7" = [double]"2008.028211216115"
    "R7" = [double]"18072.25390094504"
    "S7" = [double]"0.02116149812557349"
    "T7" = [double]"0.02116149812557349"
    "I8" = [double]"0.4193879037829277"
    "J8" = [double]"0.4193879037829278"
    "O8" = [double]"0.0001771869602491167"
    "P8" = [double]"0.0001771869602491166"
    "S8" = [double]"7.431006783654598E-05"
    "T8" = [double]"7.431006783654598E-05"
    "I9" = [double]"0.4193879037829277"
    "J9" = [double]"0.4193879037829278"
    "M9" = [double]"272.2666776666667"
    "N9" = [double]"816.800033"
    "O9" = [double]"0.4962975288350554"
    "P9" = [double]"0.4962975288350553"
    "Q9" = [double]"19750.65089529055"
    "R9" = [double]"177755.858057615"
    "S9" = [double]"0.208141180270781"
    "T9" = [double]"0.208141180270781"
    "I10" = [double]"0.4193879037829277"
    "J10" = [double]"0.4193879037829278"
    "M10" = [double]"11.73516533333333"
    "N10" = [double]"35.205496"
    "O10" = [double]"0.02139128300722342"
    "P10" = [double]"0.02139128300722341"
    "Q10" = [double]"851.2872588137467"
    "R10" = [double]"7661.585329323721"
    "S10" = [double]"0.008971245339626791"
    "T10" = [double]"0.008971245339626791"
    "I11" = [double]"0.4193879037829277"
    "J11" = [double]"0.4193879037829278"
    "M11" = [double]"236.8155566666667"
    "N11" = [double]"710.44667"
    "O11" = [double]"0.4316759456963613"
    "P11" = [double]"0.4316759456963613"
    "Q11" = [double]"17178.97109694618"
    "R11" = [double]"154610.7398725157"
    "S11" = [double]"0.1810396699791099"
    "T11" = [double]"0.1810396699791099"
    "G12" = [double]"26.10071233333333"
    "H12" = [double]"78.302137"
    "I12" = [double]"0.1508972550112184"
    "J12" = [double]"0.1508972550112184"
    "M12" = [double]"27.681071"
    "N12" = [double]"83.04321300000001"
    "O12" = [double]"0.05045805550111082"
    "P12" = [double]"0.05045805550111081"
    "Q12" = [double]"722.4956712495757"
    "R12" = [double]"6502.461041246182"
    "S12" = [double]"0.00761398206832133"
    "T12" = [double]"0.00761398206832133"
    "G13" = [double]"26.10071233333333"
    "H13" = [double]"78.302137"
    "I13" = [double]"0.1508972550112184"
    "J13" = [double]"0.1508972550112184"
    "O13" = [double]"0.0001771869602491167"
    "P13" = [double]"0.0001771869602491166"
    "Q13" = [double]"2.537093641649333"
    "R13" = [double]"22.833842774844"
    "S13" = [double]"2.673702592537357E-05"
    "T13" = [double]"2.673702592537357E-05"
    "G14" = [double]"26.10071233333333"
    "H14" = [double]"78.302137"
    "I14" = [double]"0.1508972550112184"
    "J14" = [double]"0.1508972550112184"
    "M14" = [double]"272.2666776666667"
    "N14" = [double]"816.800033"
    "O14" = [double]"0.4962975288350554"
    "P14" = [double]"0.4962975288350553"
    "Q14" = [double]"7106.354231730058"
    "R14" = [double]"63957.18808557052"
    "S14" = [double]"0.07488993477006087"
    "T14" = [double]"0.07488993477006087"
    "G15" = [double]"26.10071233333333"
    "H15" = [double]"78.302137"
    "I15" = [double]"0.1508972550112184"
    "J15" = [double]"0.1508972550112184"
    "M15" = [double]"11.73516533333333"
    "N15" = [double]"35.205496"
    "O15" = [double]"0.02139128300722342"
    "P15" = [double]"0.02139128300722341"
    "Q15" = [double]"306.2961745494391"
    "R15" = [double]"2756.665570944952"
    "S15" = [double]"0.003227885886958135"
    "T15" = [double]"0.003227885886958134"
    "G16" = [double]"26.10071233333333"
    "H16" = [double]"78.302137"
    "I16" = [double]"0.1508972550112184"
    "J16" = [double]"0.1508972550112184"
    "M16" = [double]"236.8155566666667"
    "N16" = [double]"710.44667"
    "O16" = [double]"0.4316759456963613"
    "P16" = [double]"0.4316759456963613"
    "Q16" = [double]"6181.054720614866"
    "R16" = [double]"55629.4924855338"
    "S16" = [double]"0.06513871525995268"
    "T16" = [double]"0.06513871525995268"
    "G17" = [double]"27.85999533333333"
    "H17" = [double]"83.57998600000001"
    "I17" = [double]"0.1610682791617304"
    "J17" = [double]"0.1610682791617305"
    "M17" = [double]"27.681071"
    "N17" = [double]"83.04321300000001"
    "O17" = [double]"0.05045805550111082"
    "P17" = [double]"0.05045805550111081"
    "Q17" = [double]"771.1945088816688"
    "R17" = [double]"6940.750579935019"
    "S17" = [double]"0.008127192169411006"
    "T17" = [double]"0.008127192169411006"
    "G18" = [double]"27.85999533333333"
    "H18" = [double]"83.57998600000001"
    "I18" = [double]"0.1610682791617304"
    "J18" = [double]"0.1610682791617305"
    "O18" = [double]"0.0001771869602491167"
    "P18" = [double]"0.0001771869602491166"
    "Q18" = [double]"2.708102986381333"
    "R18" = [double]"24.372926877432"
    "S18" = [double]"2.853919877722316E-05"
    "T18" = [double]"2.853919877722316E-05"
    "G19" = [double]"27.85999533333333"
    "H19" = [double]"83.57998600000001"
    "I19" = [double]"0.1610682791617304"
    "J19" = [double]"0.1610682791617305"
    "M19" = [double]"272.2666776666667"
    "N19" = [double]"816.800033"
    "O19" = [double]"0.4962975288350554"
    "P19" = [double]"0.4962975288350553"
    "Q19" = [double]"7585.348369215505"
    "R19" = [double]"68268.13532293955"
    "S19" = [double]"0.07993778892168167"
    "T19" = [double]"0.07993778892168167"
    "G20" = [double]"27.85999533333333"
    "H20" = [double]"83.57998600000001"
    "I20" = [double]"0.1610682791617304"
    "J20" = [double]"0.1610682791617305"
    "M20" = [double]"11.73516533333333"
    "N20" = [double]"35.205496"
    "O20" = [double]"0.02139128300722342"
    "P20" = [double]"0.02139128300722341"
    "Q20" = [double]"326.9416514225618"
    "R20" = [double]"2942.474862803056"
    "S20" = [double]"0.003445457143035042"
    "T20" = [double]"0.003445457143035042"
    "G21" = [double]"27.85999533333333"
    "H21" = [double]"83.57998600000001"
    "I21" = [double]"0.1610682791617304"
    "J21" = [double]"0.1610682791617305"
    "M21" = [double]"236.8155566666667"
    "N21" = [double]"710.44667"
    "O21" = [double]"0.4316759456963613"
    "P21" = [double]"0.4316759456963613"
    "Q21" = [double]"6597.680303594069"
    "R21" = [double]"59379.12273234663"
    "S21" = [double]"0.0695293017288255"
    "T21" = [double]"0.06952930172882552"
    "G22" = [double]"30.09393033333333"
    "H22" = [double]"90.281791"
    "I22" = [double]"0.1739834308659612"
    "J22" = [double]"0.1739834308659612"
    "M22" = [double]"27.681071"
    "N22" = [double]"83.04321300000001"
    "O22" = [double]"0.05045805550111082"
    "P22" = [double]"0.05045805550111081"
    "Q22" = [double]"833.0322222260537"
    "R22" = [double]"7497.290000034483"
    "S22" = [double]"0.008778865610908347"
    "T22" = [double]"0.008778865610908345"
    "G23" = [double]"30.09393033333333"
    "H23" = [double]"90.281791"
    "I23" = [double]"0.1739834308659612"
    "J23" = [double]"0.1739834308659612"
    "O23" = [double]"0.0001771869602491167"
    "P23" = [double]"0.0001771869602491166"
    "Q23" = [double]"2.925250404121333"
    "R23" = [double]"26.327253637092"
    "S23" = [double]"3.0827595248852E-05"
    "T23" = [double]"3.082759524885199E-05"
    "G24" = [double]"30.09393033333333"
    "H24" = [double]"90.281791"
    "I24" = [double]"0.1739834308659612"
    "J24" = [double]"0.1739834308659612"
    "M24" = [double]"272.2666776666667"
    "N24" = [double]"816.800033"
    "O24" = [double]"0.4962975288350554"
    "P24" = [double]"0.4962975288350553"
    "Q24" = [double]"8193.574429788789"
    "R24" = [double]"73742.1698680991"
    "S24" = [double]"0.08634754679702124"
    "T24" = [double]"0.08634754679702122"
    "G25" = [double]"30.09393033333333"
    "H25" = [double]"90.281791"
    "I25" = [double]"0.1739834308659612"
    "J25" = [double]"0.1739834308659612"
    "M25" = [double]"11.73516533333333"
    "N25" = [double]"35.205496"
    "O25" = [double]"0.02139128300722342"
    "P25" = [double]"0.02139128300722341"
    "Q25" = [double]"353.1572479914818"
    "R25" = [double]"3178.415231923336"
    "S25" = [double]"0.003721728808221466"
    "T25" = [double]"0.003721728808221465"
    "G26" = [double]"30.09393033333333"
    "H26" = [double]"90.281791"
    "I26" = [double]"0.1739834308659612"
    "J26" = [double]"0.1739834308659612"
    "M26" = [double]"236.8155566666667"
    "N26" = [double]"710.44667"
    "O26" = [double]"0.4316759456963613"
    "P26" = [double]"0.4316759456963613"
    "Q26" = [double]"7126.710864176219"
    "R26" = [double]"64140.39777758597"
    "S26" = [double]"0.07510446205456128"
    "T26" = [double]"0.07510446205456128"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
